$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1072.6182
$ws.Range("J17").Value = 1072.6182
$ws.Range("L17").Value = 3217.8546
$ws.Range("N17").Value = -3553.8546
$ws.Range("H19").Value = 809.625
$ws.Range("I19").Value = 723.875
$ws.Range("J19").Value = 895.375
$ws.Range("K19").Value = 723.875
$ws.Range("L19").Value = 895.375
$ws.Range("M19").Value = -548.875
$ws.Range("N19").Value = -1245.375
$ws.Range("H76").Value = 837317.7
$ws.Range("I76").Value = 1253587.9
$ws.Range("J76").Value = 4777.25
$ws.Range("K76").Value = 1253587.9
$ws.Range("L76").Value = 4777.25
$ws.Range("M76").Value = -1253272.9
$ws.Range("N76").Value = -5407.25
$ws.Range("H79").Value = 837317.7
$ws.Range("I79").Value = 1253587.9
$ws.Range("J79").Value = 4777.25
$ws.Range("K79").Value = 1253587.9
$ws.Range("L79").Value = 4777.25
$ws.Range("M79").Value = -1252495.9
$ws.Range("N79").Value = -6961.25
$ws.Range("H96").Value = 553.26666
$ws.Range("I96").Value = 341.1
$ws.Range("J96").Value = 977.6
$ws.Range("K96").Value = 1023.3
$ws.Range("L96").Value = 2932.8
$ws.Range("M96").Value = 349.6999999999999
$ws.Range("N96").Value = -5678.8
$ws.Range("H115").Value = 325
$ws.Range("I115").Value = 325
$ws.Range("K115").Value = 975
$ws.Range("M115").Value = 592
$ws.Range("H135").Value = 1134.075
$ws.Range("I135").Value = 977.0789
$ws.Range("K135").Value = 8793.7101
$ws.Range("M135").Value = -6258.7101
$ws.Range("H137").Value = 1756.6
$ws.Range("I137").Value = 1236.4073
$ws.Range("J137").Value = 2837
$ws.Range("K137").Value = 3709.2219
$ws.Range("L137").Value = 8511
$ws.Range("M137").Value = -1159.2219
$ws.Range("N137").Value = -13611

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7595.3335
$ws.Range("I61").Value = 995
$ws.Range("J61").Value = 100000
$ws.Range("K61").Value = 995
$ws.Range("L61").Value = 100000
$ws.Range("M61").Value = -783
$ws.Range("N61").Value = -100424
$ws.Range("H74").Value = 1324.9672
$ws.Range("I74").Value = 1150.3024
$ws.Range("K74").Value = 1150.3024
$ws.Range("M74").Value = -276.3024
$ws.Range("H77").Value = 1324.9672
$ws.Range("I77").Value = 1150.3024
$ws.Range("K77").Value = 5751.512000000001
$ws.Range("M77").Value = -1383.512000000001
$ws.Range("H132").Value = 2382.9697
$ws.Range("I132").Value = 1715.5518
$ws.Range("J132").Value = 7221.75
$ws.Range("K132").Value = 5146.6554
$ws.Range("L132").Value = 21665.25
$ws.Range("M132").Value = -2616.6554
$ws.Range("N132").Value = -26725.25
$ws.Range("H136").Value = 7595.3335
$ws.Range("I136").Value = 995
$ws.Range("J136").Value = 100000
$ws.Range("K136").Value = 2985
$ws.Range("L136").Value = 300000
$ws.Range("M136").Value = -435
$ws.Range("N136").Value = -305100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 211.5
$ws.Range("I4").Value = 207.25
$ws.Range("K4").Value = 207.25
$ws.Range("M4").Value = -92.25
$ws.Range("H94").Value = 4175953.8
$ws.Range("I94").Value = 5883440.5
$ws.Range("J94").Value = 29200
$ws.Range("K94").Value = 5883440.5
$ws.Range("L94").Value = 29200
$ws.Range("M94").Value = -5882989.5
$ws.Range("N94").Value = -30102
$ws.Range("H99").Value = 3951.3333
$ws.Range("I99").Value = 4205.2144
$ws.Range("J99").Value = 397
$ws.Range("K99").Value = 4205.2144
$ws.Range("L99").Value = 397
$ws.Range("M99").Value = -2707.2144
$ws.Range("N99").Value = -3393
$ws.Range("H134").Value = 6098.815
$ws.Range("I134").Value = 5979.2163
$ws.Range("J134").Value = 6359.1177
$ws.Range("K134").Value = 17937.6489
$ws.Range("L134").Value = 19077.3531
$ws.Range("M134").Value = -15402.6489
$ws.Range("N134").Value = -24147.3531

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 165
$ws.Range("I7").Value = 183.44444
$ws.Range("J7").Value = 131.8
$ws.Range("K7").Value = 183.44444
$ws.Range("L7").Value = 131.8
$ws.Range("M7").Value = -70.44443999999999
$ws.Range("N7").Value = -357.8
$ws.Range("H16").Value = 772.4
$ws.Range("I16").Value = 537.25
$ws.Range("K16").Value = 537.25
$ws.Range("M16").Value = -250.25
$ws.Range("H19").Value = 1875587.8
$ws.Range("I19").Value = 3000566.8
$ws.Range("J19").Value = 622.6667
$ws.Range("K19").Value = 3000566.8
$ws.Range("L19").Value = 622.6667
$ws.Range("M19").Value = -3000396.8
$ws.Range("N19").Value = -962.6667
$ws.Range("H24").Value = 1875587.8
$ws.Range("I24").Value = 3000566.8
$ws.Range("J24").Value = 622.6667
$ws.Range("K24").Value = 3000566.8
$ws.Range("L24").Value = 622.6667
$ws.Range("M24").Value = -3000396.8
$ws.Range("N24").Value = -962.6667
$ws.Range("H58").Value = 1532.421
$ws.Range("I58").Value = 1652
$ws.Range("J58").Value = 894.6667
$ws.Range("K58").Value = 1652
$ws.Range("L58").Value = 894.6667
$ws.Range("M58").Value = -1449
$ws.Range("N58").Value = -1300.6667
$ws.Range("H113").Value = 772.4
$ws.Range("I113").Value = 537.25
$ws.Range("K113").Value = 537.25
$ws.Range("M113").Value = 1632.75
$ws.Range("H118").Value = 215000
$ws.Range("J118").Value = 215000
$ws.Range("L118").Value = 215000
$ws.Range("N118").Value = -218314
$ws.Range("H132").Value = 3402.138
$ws.Range("I132").Value = 3252.9473
$ws.Range("J132").Value = 3685.6
$ws.Range("K132").Value = 9758.841899999999
$ws.Range("L132").Value = 11056.8
$ws.Range("M132").Value = -7228.841899999999
$ws.Range("N132").Value = -16116.8
$ws.Range("H136").Value = 1532.421
$ws.Range("I136").Value = 1652
$ws.Range("J136").Value = 894.6667
$ws.Range("K136").Value = 4956
$ws.Range("L136").Value = 2684.0001
$ws.Range("M136").Value = -2406
$ws.Range("N136").Value = -7784.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 461.46667
$ws.Range("I8").Value = 461.46667
$ws.Range("K8").Value = 1384.40001
$ws.Range("M8").Value = -1245.40001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1616.4286
$ws.Range("I132").Value = 1616.4286
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4849.2858
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2319.2858
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 39100.08
$ws.Range("J136").Value = 39100.08
$ws.Range("L136").Value = 117300.24
$ws.Range("N136").Value = -122400.24

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1999.1666
$ws.Range("I61").Value = 1999.1666
$ws.Range("K61").Value = 1999.1666
$ws.Range("M61").Value = -1797.1666
$ws.Range("H100").Value = 1677.8
$ws.Range("I100").Value = 1597.25
$ws.Range("K100").Value = 1597.25
$ws.Range("M100").Value = -1056.25
$ws.Range("H113").Value = 1999.1666
$ws.Range("I113").Value = 1999.1666
$ws.Range("K113").Value = 1999.1666
$ws.Range("M113").Value = 170.8334
$ws.Range("H132").Value = 1999.7073
$ws.Range("I132").Value = 1157.6666
$ws.Range("K132").Value = 3472.9998
$ws.Range("M132").Value = -942.9998000000001
$ws.Range("H136").Value = 1567.2941
$ws.Range("I136").Value = 1001.7308
$ws.Range("K136").Value = 3005.1924
$ws.Range("M136").Value = -455.1923999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 721.4
$ws.Range("I113").Value = 481.26315
$ws.Range("K113").Value = 1443.78945
$ws.Range("M113").Value = 726.21055
$ws.Range("H132").Value = 2195.5435
$ws.Range("I132").Value = 2062.3076
$ws.Range("J132").Value = 2937.8572
$ws.Range("K132").Value = 6186.9228
$ws.Range("L132").Value = 8813.571599999999
$ws.Range("M132").Value = -3656.9228
$ws.Range("N132").Value = -13873.5716
$ws.Range("H136").Value = 11666.5
$ws.Range("I136").Value = 12750
$ws.Range("J136").Value = 9499.5
$ws.Range("K136").Value = 38250
$ws.Range("L136").Value = 28498.5
$ws.Range("M136").Value = -35700
$ws.Range("N136").Value = -33598.5
